$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 44533
$ws.Range("L2").Value2 = 'Primera'
$ws.Range("M2").Value2 = 140
$ws.Range("N2").Value2 = 14000
$ws.Range("O2").Value2 = 15000
$ws.Range("P2").Value2 = 14500
$ws.Range("Q2").Value2 = '$/caja 10 kilos'
$ws.Range("R2").Value2 = 'Región de O''Higgins'
$ws.Range("S2").Value2 = 1450
$ws.Range("T2").Value2 = 10

$ws.Range("D3").Value2 = 44895
$ws.Range("L3").Value2 = 'Segunda'
$ws.Range("M3").Value2 = 130
$ws.Range("N3").Value2 = 19000
$ws.Range("O3").Value2 = 20000
$ws.Range("P3").Value2 = 19462
$ws.Range("Q3").Value2 = '$/caja 16 kilos granel'
$ws.Range("R3").Value2 = 'Región de O''Higgins'
$ws.Range("S3").Value2 = 1216
$ws.Range("T3").Value2 = 16

$ws.Range("D4").Value2 = 44174
$ws.Range("L4").Value2 = 'Primera'
$ws.Range("M4").Value2 = 300
$ws.Range("N4").Value2 = 19000
$ws.Range("O4").Value2 = 20000
$ws.Range("P4").Value2 = 19500
$ws.Range("Q4").Value2 = '$/bandeja 18 kilos'
$ws.Range("R4").Value2 = 'Región Metropolitana'
$ws.Range("S4").Value2 = 1083
$ws.Range("T4").Value2 = 18

$ws.Range("D5").Value2 = 44880
$ws.Range("L5").Value2 = 'Primera'
$ws.Range("M5").Value2 = 200
$ws.Range("N5").Value2 = 33000
$ws.Range("O5").Value2 = 34000
$ws.Range("P5").Value2 = 33500
$ws.Range("Q5").Value2 = '$/caja 10 kilos'
$ws.Range("R5").Value2 = 'Región de O''Higgins'
$ws.Range("S5").Value2 = 3350
$ws.Range("T5").Value2 = 10

$ws.Range("D6").Value2 = 44169
$ws.Range("L6").Value2 = 'Primera'
$ws.Range("M6").Value2 = 250
$ws.Range("N6").Value2 = 20000
$ws.Range("O6").Value2 = 22000
$ws.Range("P6").Value2 = 21000
$ws.Range("Q6").Value2 = '$/bandeja 18 kilos'
$ws.Range("R6").Value2 = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S6").Value2 = 1167
$ws.Range("T6").Value2 = 18

$ws.Range("D7").Value2 = 44881
$ws.Range("L7").Value2 = 'Segunda'
$ws.Range("M7").Value2 = 300
$ws.Range("N7").Value2 = 41000
$ws.Range("O7").Value2 = 42000
$ws.Range("P7").Value2 = 41500
$ws.Range("Q7").Value2 = '$/bandeja 18 kilos'
$ws.Range("R7").Value2 = 'Región de Coquimbo'
$ws.Range("S7").Value2 = 2306
$ws.Range("T7").Value2 = 18

$ws.Range("D8").Value2 = 44545
$ws.Range("L8").Value2 = 'Primera'
$ws.Range("M8").Value2 = 200
$ws.Range("N8").Value2 = 24000
$ws.Range("O8").Value2 = 25000
$ws.Range("P8").Value2 = 24500
$ws.Range("Q8").Value2 = '$/bandeja 18 kilos'
$ws.Range("R8").Value2 = 'Región de Coquimbo'
$ws.Range("S8").Value2 = 1361
$ws.Range("T8").Value2 = 18

$ws.Range("D9").Value2 = 44894
$ws.Range("L9").Value2 = 'Segunda'
$ws.Range("M9").Value2 = 130
$ws.Range("N9").Value2 = 19000
$ws.Range("O9").Value2 = 20000
$ws.Range("P9").Value2 = 19462
$ws.Range("Q9").Value2 = '$/caja 16 kilos granel'
$ws.Range("R9").Value2 = 'Región de O''Higgins'
$ws.Range("S9").Value2 = 1216
$ws.Range("T9").Value2 = 16

$ws.Range("D10").Value2 = 44160
$ws.Range("L10").Value2 = 'Primera'
$ws.Range("M10").Value2 = 250
$ws.Range("N10").Value2 = 24000
$ws.Range("O10").Value2 = 25000
$ws.Range("P10").Value2 = 24500
$ws.Range("Q10").Value2 = '$/bandeja 18 kilos'
$ws.Range("R10").Value2 = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S10").Value2 = 1361
$ws.Range("T10").Value2 = 18

$ws.Range("D11").Value2 = 44901
$ws.Range("L11").Value2 = 'Segunda'
$ws.Range("M11").Value2 = 200
$ws.Range("N11").Value2 = 17000
$ws.Range("O11").Value2 = 18000
$ws.Range("P11").Value2 = 17500
$ws.Range("Q11").Value2 = '$/bandeja 18 kilos'
$ws.Range("R11").Value2 = 'Región de O''Higgins'
$ws.Range("S11").Value2 = 972
$ws.Range("T11").Value2 = 18

$ws.Range("D13").Value2 = 44917
$ws.Range("L13").Value2 = 'Segunda'
$ws.Range("M13").Value2 = 250
$ws.Range("N13").Value2 = 20000
$ws.Range("O13").Value2 = 23000
$ws.Range("P13").Value2 = 21800
$ws.Range("Q13").Value2 = '$/caja 18 kilos'
$ws.Range("R13").Value2 = 'Región de Coquimbo'
$ws.Range("S13").Value2 = 1211
$ws.Range("T13").Value2 = 18

$ws.Range("D14").Value2 = 44524
$ws.Range("L14").Value2 = 'Segunda'
$ws.Range("M14").Value2 = 200
$ws.Range("N14").Value2 = 27000
$ws.Range("O14").Value2 = 28000
$ws.Range("P14").Value2 = 27500
$ws.Range("Q14").Value2 = '$/bandeja 18 kilos'
$ws.Range("R14").Value2 = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S14").Value2 = 1528
$ws.Range("T14").Value2 = 18
